$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1311
$ws.Range("F3").Value = 320.7
$ws.Range("F4").Value = 78.52
$ws.Range("F5").Value = 101.5
$ws.Range("F6").Value = 0.008772
$ws.Range("F7").Value = 6.235
$ws.Range("F8").Value = 5.625

$ws.Range("F9").Select()
